# Update the "想去人数" (attendance count) figures in column F
# for the 展览 (sheet 1) and 全部类型 (sheet 4) worksheets, to match a
# newer data pull.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 228
    "F3"  = 439
    "F4"  = 13033
    "F5"  = 1342
    "F6"  = 219
    "F9"  = 166
    "F10" = 228
    "F11" = 470
    "F13" = 68
    "F17" = 415
    "F18" = 5543
    "F19" = 108
    "F20" = 56
    "F21" = 961
    "F22" = 34
    "F24" = 144
}

# Sheet 1 = "展览" and Sheet 4 = "全部类型" both carry the same rows of
# event data and both received the updated counts.
$sheetIndexes = @(1, 4)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
